# ADD results from server
# Update computed result values (row 2) on each year sheet with fresh
# server-provided figures. Columns C, D, F, H, J, K remain 0 (unchanged).

$wb = $excel.ActiveWorkbook

$data = @{
    "2025" = @{ A=0;                 B=277.8902526399997;  E=29092.72506141524; G=8095.925712661834; I=14865.25872276;    L=50912.59821312752; M=11247.09127927;    N=7270.39941619107;  O=6890.515200515631 }
    "2030" = @{ A=0;                 B=5522.228665642542;  E=56035.50691766649; G=8095.925712661834; I=37516.80488821431; L=72537.0054129926;  M=21815.58876574551; N=10913.44862569963; O=9427.219998659948 }
    "2035" = @{ A=2266.48797166071;  B=7059.768887134545;  E=67081.2042405315;  G=8095.925712661834; I=53800.46711853385; L=72537.0054129926;  M=27679.0502989715;  N=15939.31202052396; O=15294.41615167074 }
    "2040" = @{ A=2266.48797166071;  B=7059.768887134545;  E=67081.2042405315;  G=8095.925712661834; I=53800.46711853385; L=72537.0054129926;  M=27679.0502989715;  N=15939.31202052396; O=15294.41615167074 }
    "2045" = @{ A=2266.48797166071;  B=7059.768887134545;  E=67081.2042405315;  G=8095.925712661834; I=53800.46711853385; L=72537.0054129926;  M=27679.0502989715;  N=15939.31202052396; O=15294.41615167074 }
    "2050" = @{ A=2266.48797166071;  B=7059.768887134545;  E=67081.2042405315;  G=8095.925712661834; I=53800.46711853385; L=72537.0054129926;  M=27679.0502989715;  N=15939.31202052396; O=15294.41615167074 }
}

foreach ($sheetName in $data.Keys) {
    $name = [string]$sheetName
    $ws = $wb.Worksheets($name)
    $vals = $data[$sheetName]
    foreach ($col in $vals.Keys) {
        $colName = [string]$col
        $ws.Range("$colName`2").Value = $vals[$col]
    }
}
